{"js": "// Rewrite the report body to match the target revision:\n//  - split the author/byline paragraph into separate runs (with a right\n//    tab stop at 8931 twips) and add a second author \"Kunga Kartung \u2013 s3588773\"\n//  - add \"Analysis\" and \"Solution\" heading paragraphs (bold, 28 half-points)\n//    separated from neighbours by blank paragraphs\n//  - rewrite / extend the analysis paragraphs and append several new\n//    paragraphs that make up the new \"Solution\" + conclusion section\n//\n// The cleanest way to reproduce this very large, highly granular set of\n// run-splits/insertions exactly (including <w:proofErr>, <w:tab/>,\n// the <w:bookmarkStart/End> around \"_GoBack\" and xml:space=\"preserve\")\n// is to replace the whole body with the target OOXML in one shot.\nconst body = context.document.body;\nconst ooxml = \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:pPr><w:pStyle w:val=\\\"Title\\\"/></w:pPr><w:r><w:t>Assignment 1 \\u2013 Report</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val=\\\"right\\\" w:pos=\\\"8931\\\"/></w:tabs><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">Aleksandar </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t>Stefanovic</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> \\u2013 s3</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t>605170</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:tab/></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t>Kunga</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t>Kartung</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> \\u2013 s3588773</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t>Analysis</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">The site is </w:t></w:r><w:r><w:t>clear in layout</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> and has all the required information, however the layout of content could be reformed, reducing the navigation needed to reach the relevant pages.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t>T</w:t></w:r><w:r><w:t>he tables used for presenting past and in progress tickets are in contrast with the simplicity of the rest of the site, being heavy with detail, a weak point in the sites design.</w:t></w:r></w:p><w:p><w:r><w:t>The landing page for the site is clear in design,</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> apart from a shopping cart being used as an icon to request a service</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">. Even so, the </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">overall </w:t></w:r><w:r><w:t>function of the page isn\\u2019t clear, with little information being conveyed that isn\\u2019t already found in the navigation bar. With that in mind, the page could be improved by adding information to explain how to use the site and provide context or the page could be removed, assuming the layout of the site is self-explanatory.</w:t></w:r></w:p><w:p><w:r><w:t>The \\u2018request service\\u2019</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> page lists services available for various areas in a tiled series of boxes. Even though the design isn\\u2019t bad, it could be simplified into just a list, with subheadi</w:t></w:r><w:r><w:t>ngs replacing the various tiles, allowing for a more mobile friendly layout and not confusing the user with the panel-like divisions that may seem like buttons.</w:t></w:r></w:p><w:p><w:r><w:t>The ticketing</w:t></w:r><w:r><w:t>/report an issue</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> page is set up concisely, </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">each form input being clear and easy to interact with, except for the attachments icon out of sight in the top right of the form. </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">The dropdown inputs have a default </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">\\u201c-- </w:t></w:r><w:r><w:t>None</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> --\\u201d</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> option that shows up when selecting another option, rather than having the default as just a placeholder. </w:t></w:r><w:r><w:t>When submitting a form, an alert will show up when inputs are missing, unnecessary given that missing or invalid inputs are given a clear red outline. The ability to collapse parts of the form is another feature that could be removed, as its utility is limited in such a short form. Also, the basic bootstrap styling could be redone to convey a style closer the rest of RMIT\\u2019s pages.</w:t></w:r></w:p><w:p><w:r><w:t>The \\u2018user card\\u2019, as can be seen by hovering over the \\u2018</w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t>\\u2019 icon next to the name in the ticketing pages, is another detailed feature that could be improved by replacing it. Currently, to edit the card you need to hold shift and then the window persists. This is not only cumbersome, but also requires the use of a keyboard, and it is not evident how a mobile device or other pointing devices will be able to interface with the card. To replace this, either a separate page can be used to capture this data, allowing for all devices to change it in a simple series of form inputs. On the other hand, this feature could be removed altogether depending on its utility to the user.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">The </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">following </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">track progress and closed enquiries pages </w:t></w:r><w:r><w:t>share</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> the same weaknesses, with good design barred by a heavy and overcomplicated table system. The customization and search capabilities of the system are outside of the needs of a user that will not have enough tickets to support search and filters at the expense of loading times and a</w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t>infle</w:t></w:r><w:r><w:t>xible layout that will not translate well in</w:t></w:r><w:r><w:t>to mobile.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> The utility of separating the tracked and closed enquires is also debatable, especially if the large tables were compacted to a smaller design by just dis</w:t></w:r><w:r><w:t>playing the tickets as separate \\u2018cards\\u2019.</w:t></w:r></w:p><w:p><w:r><w:t>Whe</w:t></w:r><w:r><w:t>n requesting a general IT enquiry from the \\u2018Request a Service\\u2019 page</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">, the layout </w:t></w:r><w:r><w:t>has similar</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> problems as the ticketing page </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">from </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">before. The same collapsing function, attachments icon and </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space=\\\"preserve\\\">obtrusive alert can be found, with their corresponding drawbacks. In addition, the submit button on the side may also pose a problem, </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">as </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">the scrolling animation being used may </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">render the </w:t></w:r><w:r><w:t>form</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> unusable</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> depending on</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> the size of the screen. Keeping the button fixed at the bottom of the form does not subtract from the user\\u2019s experience, while being more reliable.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t>Solution</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">To help simplify the site, all the superfluous elements identified before should be removed, while retaining the necessary inputs for submitting tickets and </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">viewing them. </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">Firstly, the home page can be removed from the site, there is little information conveyed in its current form, and the site is simple enough that the </w:t></w:r><w:r><w:t>user can navigate and understand through just the scroll bar</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Creating a ticket can now</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> be</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> streamlined, with the \\u2018Request a Service\\u2019 and \\u2018Report Issue\\u2019 pages being merged into one and the listing of services not being required. </w:t></w:r><w:r><w:t>The form fields used to report an issue can be preserved, just changing the \\u2018customer details\\u2019 section to compensate for the lack of login. Mobile-friendly scaling can also be added to size the forms appropriately.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">Viewing completed and in progress tickets can similarly be simplified, with the removal of the table system for </w:t></w:r><w:r><w:t>both</w:t></w:r><w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/><w:bookmarkEnd w:id=\\\"0\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> types of tickets and removing the further separation found in the three tables on both the \\u2018Track Progress\\u2019 and \\u2018Closed Enquiries\\u2019 pages on the current support site. Instead, the two pages will be merged into one, with subheadings separating </w:t></w:r><w:r><w:t>the completed and in progress tickets.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> In addition, each ticket will be </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">put into a \\u2018card\\u2019 for layout, allowing </w:t></w:r><w:r><w:t>for an intuitive</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> and scalable des</w:t></w:r><w:r><w:t>ign that has the ticket, post information and comments all in one place. However, this design removes the granular control of the table, sorting could be implemented but limited. This compromise aligns the site with the needs of the user, eliminating advanced features that do more harm than good.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">In conclusion, the largest problems with the site arise from overly granular, complex controls and unnecessary content that </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">can be removed, merged or simplified. In doing this, we create a simpler interface that still </w:t></w:r><w:r><w:t>satisfies the user\\u2019s needs, and can scale to smaller screens easily.</w:t></w:r></w:p><w:sectPr><w:pgSz w:w=\\\"11906\\\" w:h=\\\"16838\\\"/><w:pgMar w:top=\\\"1440\\\" w:right=\\\"1440\\\" w:bottom=\\\"1440\\\" w:left=\\\"1440\\\" w:header=\\\"708\\\" w:footer=\\\"708\\\" w:gutter=\\\"0\\\"/><w:cols w:space=\\\"708\\\"/><w:docGrid w:linePitch=\\\"360\\\"/></w:sectPr>\\n</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\nbody.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Rewrite the report body to match the target revision:\n#  - split the author/byline paragraph into separate runs (with a right\n#    tab stop at 8931 twips) and add a second author \"Kunga Kartung - s3588773\"\n#  - add \"Analysis\" and \"Solution\" heading paragraphs (bold, 28 half-points)\n#    separated from neighbours by blank paragraphs\n#  - rewrite / extend the analysis paragraphs and append several new\n#    paragraphs that make up the new \"Solution\" + conclusion section\n#\n# The cleanest way to reproduce this very large, highly granular set of\n# run-splits/insertions exactly (including proofErr markers, manual tabs,\n# the bookmarkStart/End around \"_GoBack\" and preserved whitespace runs)\n# is to replace the whole document body with the target OOXML in one shot.\n$d = $word.ActiveDocument\n$xml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"Title\"/></w:pPr><w:r><w:t>Assignment 1 \u2013 Report</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val=\"right\" w:pos=\"8931\"/></w:tabs><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\">Aleksandar </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Stefanovic</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\"> \u2013 s3</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>605170</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:tab/></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Kunga</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Kartung</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\"> \u2013 s3588773</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Analysis</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">The site is </w:t></w:r><w:r><w:t>clear in layout</w:t></w:r><w:r><w:t xml:space=\"preserve\"> and has all the required information, however the layout of content could be reformed, reducing the navigation needed to reach the relevant pages.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>T</w:t></w:r><w:r><w:t>he tables used for presenting past and in progress tickets are in contrast with the simplicity of the rest of the site, being heavy with detail, a weak point in the sites design.</w:t></w:r></w:p><w:p><w:r><w:t>The landing page for the site is clear in design,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> apart from a shopping cart being used as an icon to request a service</w:t></w:r><w:r><w:t xml:space=\"preserve\">. Even so, the </w:t></w:r><w:r><w:t xml:space=\"preserve\">overall </w:t></w:r><w:r><w:t>function of the page isn\u2019t clear, with little information being conveyed that isn\u2019t already found in the navigation bar. With that in mind, the page could be improved by adding information to explain how to use the site and provide context or the page could be removed, assuming the layout of the site is self-explanatory.</w:t></w:r></w:p><w:p><w:r><w:t>The \u2018request service\u2019</w:t></w:r><w:r><w:t xml:space=\"preserve\"> page lists services available for various areas in a tiled series of boxes. Even though the design isn\u2019t bad, it could be simplified into just a list, with subheadi</w:t></w:r><w:r><w:t>ngs replacing the various tiles, allowing for a more mobile friendly layout and not confusing the user with the panel-like divisions that may seem like buttons.</w:t></w:r></w:p><w:p><w:r><w:t>The ticketing</w:t></w:r><w:r><w:t>/report an issue</w:t></w:r><w:r><w:t xml:space=\"preserve\"> page is set up concisely, </w:t></w:r><w:r><w:t xml:space=\"preserve\">each form input being clear and easy to interact with, except for the attachments icon out of sight in the top right of the form. </w:t></w:r><w:r><w:t xml:space=\"preserve\">The dropdown inputs have a default </w:t></w:r><w:r><w:t xml:space=\"preserve\">\u201c-- </w:t></w:r><w:r><w:t>None</w:t></w:r><w:r><w:t xml:space=\"preserve\"> --\u201d</w:t></w:r><w:r><w:t xml:space=\"preserve\"> option that shows up when selecting another option, rather than having the default as just a placeholder. </w:t></w:r><w:r><w:t>When submitting a form, an alert will show up when inputs are missing, unnecessary given that missing or invalid inputs are given a clear red outline. The ability to collapse parts of the form is another feature that could be removed, as its utility is limited in such a short form. Also, the basic bootstrap styling could be redone to convey a style closer the rest of RMIT\u2019s pages.</w:t></w:r></w:p><w:p><w:r><w:t>The \u2018user card\u2019, as can be seen by hovering over the \u2018</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>\u2019 icon next to the name in the ticketing pages, is another detailed feature that could be improved by replacing it. Currently, to edit the card you need to hold shift and then the window persists. This is not only cumbersome, but also requires the use of a keyboard, and it is not evident how a mobile device or other pointing devices will be able to interface with the card. To replace this, either a separate page can be used to capture this data, allowing for all devices to change it in a simple series of form inputs. On the other hand, this feature could be removed altogether depending on its utility to the user.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">The </w:t></w:r><w:r><w:t xml:space=\"preserve\">following </w:t></w:r><w:r><w:t xml:space=\"preserve\">track progress and closed enquiries pages </w:t></w:r><w:r><w:t>share</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the same weaknesses, with good design barred by a heavy and overcomplicated table system. The customization and search capabilities of the system are outside of the needs of a user that will not have enough tickets to support search and filters at the expense of loading times and a</w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>infle</w:t></w:r><w:r><w:t>xible layout that will not translate well in</w:t></w:r><w:r><w:t>to mobile.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> The utility of separating the tracked and closed enquires is also debatable, especially if the large tables were compacted to a smaller design by just dis</w:t></w:r><w:r><w:t>playing the tickets as separate \u2018cards\u2019.</w:t></w:r></w:p><w:p><w:r><w:t>Whe</w:t></w:r><w:r><w:t>n requesting a general IT enquiry from the \u2018Request a Service\u2019 page</w:t></w:r><w:r><w:t xml:space=\"preserve\">, the layout </w:t></w:r><w:r><w:t>has similar</w:t></w:r><w:r><w:t xml:space=\"preserve\"> problems as the ticketing page </w:t></w:r><w:r><w:t xml:space=\"preserve\">from </w:t></w:r><w:r><w:t xml:space=\"preserve\">before. The same collapsing function, attachments icon and </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">obtrusive alert can be found, with their corresponding drawbacks. In addition, the submit button on the side may also pose a problem, </w:t></w:r><w:r><w:t xml:space=\"preserve\">as </w:t></w:r><w:r><w:t xml:space=\"preserve\">the scrolling animation being used may </w:t></w:r><w:r><w:t xml:space=\"preserve\">render the </w:t></w:r><w:r><w:t>form</w:t></w:r><w:r><w:t xml:space=\"preserve\"> unusable</w:t></w:r><w:r><w:t xml:space=\"preserve\"> depending on</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the size of the screen. Keeping the button fixed at the bottom of the form does not subtract from the user\u2019s experience, while being more reliable.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Solution</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">To help simplify the site, all the superfluous elements identified before should be removed, while retaining the necessary inputs for submitting tickets and </w:t></w:r><w:r><w:t xml:space=\"preserve\">viewing them. </w:t></w:r><w:r><w:t xml:space=\"preserve\">Firstly, the home page can be removed from the site, there is little information conveyed in its current form, and the site is simple enough that the </w:t></w:r><w:r><w:t>user can navigate and understand through just the scroll bar</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Creating a ticket can now</w:t></w:r><w:r><w:t xml:space=\"preserve\"> be</w:t></w:r><w:r><w:t xml:space=\"preserve\"> streamlined, with the \u2018Request a Service\u2019 and \u2018Report Issue\u2019 pages being merged into one and the listing of services not being required. </w:t></w:r><w:r><w:t>The form fields used to report an issue can be preserved, just changing the \u2018customer details\u2019 section to compensate for the lack of login. Mobile-friendly scaling can also be added to size the forms appropriately.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">Viewing completed and in progress tickets can similarly be simplified, with the removal of the table system for </w:t></w:r><w:r><w:t>both</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\"> types of tickets and removing the further separation found in the three tables on both the \u2018Track Progress\u2019 and \u2018Closed Enquiries\u2019 pages on the current support site. Instead, the two pages will be merged into one, with subheadings separating </w:t></w:r><w:r><w:t>the completed and in progress tickets.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> In addition, each ticket will be </w:t></w:r><w:r><w:t xml:space=\"preserve\">put into a \u2018card\u2019 for layout, allowing </w:t></w:r><w:r><w:t>for an intuitive</w:t></w:r><w:r><w:t xml:space=\"preserve\"> and scalable des</w:t></w:r><w:r><w:t>ign that has the ticket, post information and comments all in one place. However, this design removes the granular control of the table, sorting could be implemented but limited. This compromise aligns the site with the needs of the user, eliminating advanced features that do more harm than good.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">In conclusion, the largest problems with the site arise from overly granular, complex controls and unnecessary content that </w:t></w:r><w:r><w:t xml:space=\"preserve\">can be removed, merged or simplified. In doing this, we create a simpler interface that still </w:t></w:r><w:r><w:t>satisfies the user\u2019s needs, and can scale to smaller screens easily.</w:t></w:r></w:p><w:sectPr><w:pgSz w:w=\"11906\" w:h=\"16838\"/><w:pgMar w:top=\"1440\" w:right=\"1440\" w:bottom=\"1440\" w:left=\"1440\" w:header=\"708\" w:footer=\"708\" w:gutter=\"0\"/><w:cols w:space=\"708\"/><w:docGrid w:linePitch=\"360\"/></w:sectPr>\n</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$d.Content.InsertXML($xml)\n"}
